# Update stack trace line numbers in the error text embedded in the
# document, reflecting the source line shifts caused by moving the
# M2Doc project version from 2.0.0 to 2.0.1.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# M2DocEvaluator.java line renumbering
Replace-Text "M2DocEvaluator.java:543)" "M2DocEvaluator.java:555)"
Replace-Text "M2DocEvaluator.java:1300)" "M2DocEvaluator.java:1305)"
Replace-Text "M2DocEvaluator.java:278)" "M2DocEvaluator.java:283)"
Replace-Text "M2DocEvaluator.java:267)" "M2DocEvaluator.java:272)"

# M2DocEvaluator.java:1084 appears three times in the stack trace;
# wdReplaceAll (last argument = 2) replaces every occurrence in one call.
Replace-Text "M2DocEvaluator.java:1084)" "M2DocEvaluator.java:1096)"

# AbstractTemplatesTestSuite.java line renumbering
Replace-Text "AbstractTemplatesTestSuite.java:475)" "AbstractTemplatesTestSuite.java:479)"
Replace-Text "AbstractTemplatesTestSuite.java:384)" "AbstractTemplatesTestSuite.java:388)"
